$d = $word.ActiveDocument

# The document has two logos, each embedded twice (once in the "first page"
# header/footer, once in the "default" header/footer):
#   - "BTec_Logo-Orange"  (media/image1.jpg)  -> rename the picture to image2.jpg
#   - "...PearsonLogo.png" (media/image2.png) -> rename the picture to image1.png
# Walk every header/footer of every section and rename the inline picture(s)
# found there based on their (unchanged) alt-text/description, so the script
# is robust regardless of which physical header/footer part Word exposes as
# index 1 vs 2.

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections($s)

    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers($h)
        if ($hdr.Exists) {
            $ishapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $ishapes.Count; $i++) {
                $shp = $ishapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers($f)
        if ($ftr.Exists) {
            $ishapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $ishapes.Count; $i++) {
                $shp = $ishapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
